$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H28").Value = 434
$ws.Range("I28").Value = 426.66666
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 426.66666
$ws.Range("L28").Value = 500
$ws.Range("M28").Value = 58.33334000000002
$ws.Range("N28").Value = -1470
$ws.Range("H40").Value = 4057
$ws.Range("J40").Value = 2883.1667
$ws.Range("L40").Value = 2883.1667
$ws.Range("N40").Value = -3233.1667
$ws.Range("H107").Value = 735.9286
$ws.Range("I107").Value = 472.14285
$ws.Range("J107").Value = 999.7143
$ws.Range("K107").Value = 472.14285
$ws.Range("L107").Value = 999.7143
$ws.Range("M107").Value = 1447.85715
$ws.Range("N107").Value = -4839.7143
$ws.Range("H113").Value = 2274.318
$ws.Range("I113").Value = 1721.3636
$ws.Range("J113").Value = 2827.2727
$ws.Range("K113").Value = 1721.3636
$ws.Range("L113").Value = 2827.2727
$ws.Range("M113").Value = 1532.6364
$ws.Range("N113").Value = -9335.2727
$ws.Range("H129").Value = 6581.921
$ws.Range("J129").Value = 8231.233
$ws.Range("L129").Value = 24693.699
$ws.Range("N129").Value = -34693.699
$ws.Range("H138").Value = 3305.9678
$ws.Range("I138").Value = 2708.5
$ws.Range("J138").Value = 3420.8655
$ws.Range("K138").Value = 8125.5
$ws.Range("L138").Value = 10262.5965
$ws.Range("M138").Value = -2985.5
$ws.Range("N138").Value = -20542.5965

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1093.9656
$ws.Range("I2").Value = 966.8333
$ws.Range("K2").Value = 966.8333
$ws.Range("M2").Value = -853.8333
$ws.Range("H45").Value = 2223.75
$ws.Range("I45").Value = 2398.125
$ws.Range("J45").Value = 1875
$ws.Range("K45").Value = 2398.125
$ws.Range("L45").Value = 1875
$ws.Range("M45").Value = -2021.125
$ws.Range("N45").Value = -2629
$ws.Range("H61").Value = 2091.28
$ws.Range("I61").Value = 1447
$ws.Range("J61").Value = 3236.6667
$ws.Range("K61").Value = 1447
$ws.Range("L61").Value = 3236.6667
$ws.Range("M61").Value = -1235
$ws.Range("N61").Value = -3660.6667
$ws.Range("H74").Value = 2488.9546
$ws.Range("I74").Value = 2395.8708
$ws.Range("J74").Value = 2710.923
$ws.Range("K74").Value = 2395.8708
$ws.Range("L74").Value = 2710.923
$ws.Range("M74").Value = -1521.8708
$ws.Range("N74").Value = -4458.923
$ws.Range("H77").Value = 2488.9546
$ws.Range("I77").Value = 2395.8708
$ws.Range("J77").Value = 2710.923
$ws.Range("K77").Value = 11979.354
$ws.Range("L77").Value = 13554.615
$ws.Range("M77").Value = -7611.354000000001
$ws.Range("N77").Value = -22290.615
$ws.Range("H116").Value = 1093.9656
$ws.Range("I116").Value = 966.8333
$ws.Range("K116").Value = 966.8333
$ws.Range("M116").Value = 1327.1667
$ws.Range("H122").Value = 2072.5
$ws.Range("I122").Value = 1745.8182
$ws.Range("J122").Value = 3612.5715
$ws.Range("K122").Value = 5237.4546
$ws.Range("L122").Value = 10837.7145
$ws.Range("M122").Value = -2787.4546
$ws.Range("N122").Value = -15737.7145
$ws.Range("H136").Value = 2091.28
$ws.Range("I136").Value = 1447
$ws.Range("J136").Value = 3236.6667
$ws.Range("K136").Value = 4341
$ws.Range("L136").Value = 9710.000100000001
$ws.Range("M136").Value = -1791
$ws.Range("N136").Value = -14810.0001

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1093.9656
$ws.Range("I3").Value = 966.8333
$ws.Range("K3").Value = 966.8333
$ws.Range("M3").Value = -852.8333
$ws.Range("H86").Value = 6251787
$ws.Range("I86").Value = 7409084.5
$ws.Range("J86").Value = 2379.8
$ws.Range("K86").Value = 7409084.5
$ws.Range("L86").Value = 2379.8
$ws.Range("M86").Value = -7407961.5
$ws.Range("N86").Value = -4625.8
$ws.Range("H89").Value = 6251787
$ws.Range("I89").Value = 7409084.5
$ws.Range("J89").Value = 2379.8
$ws.Range("K89").Value = 37045422.5
$ws.Range("L89").Value = 11899
$ws.Range("M89").Value = -37039806.5
$ws.Range("N89").Value = -23131
$ws.Range("H105").Value = 1797.6666
$ws.Range("I105").Value = 1160.625
$ws.Range("J105").Value = 2525.7144
$ws.Range("K105").Value = 1160.625
$ws.Range("L105").Value = 2525.7144
$ws.Range("M105").Value = 586.375
$ws.Range("N105").Value = -6019.7144
$ws.Range("H107").Value = 1607.5
$ws.Range("I107").Value = 1479.5834
$ws.Range("K107").Value = 1479.5834
$ws.Range("M107").Value = 440.4166

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 33334152
$ws.Range("I16").Value = 55556164
$ws.Range("J16").Value = 1135.5
$ws.Range("K16").Value = 55556164
$ws.Range("L16").Value = 1135.5
$ws.Range("M16").Value = -55555877
$ws.Range("N16").Value = -1709.5
$ws.Range("H99").Value = 5959766
$ws.Range("I99").Value = 7150119
$ws.Range("J99").Value = 8000
$ws.Range("K99").Value = 7150119
$ws.Range("L99").Value = 8000
$ws.Range("M99").Value = -7148621
$ws.Range("N99").Value = -10996
$ws.Range("H107").Value = 71431090
$ws.Range("I107").Value = 100002170
$ws.Range("J107").Value = 3400
$ws.Range("K107").Value = 100002170
$ws.Range("L107").Value = 3400
$ws.Range("M107").Value = -100000250
$ws.Range("N107").Value = -7240
$ws.Range("H113").Value = 33334152
$ws.Range("I113").Value = 55556164
$ws.Range("J113").Value = 1135.5
$ws.Range("K113").Value = 55556164
$ws.Range("L113").Value = 1135.5
$ws.Range("M113").Value = -55553994
$ws.Range("N113").Value = -5475.5
$ws.Range("H122").Value = 5129212
$ws.Range("I122").Value = 8334005.5
$ws.Range("J122").Value = 1542
$ws.Range("K122").Value = 25002016.5
$ws.Range("L122").Value = 4626
$ws.Range("M122").Value = -24999566.5
$ws.Range("N122").Value = -9526
$ws.Range("H126").Value = 5959766
$ws.Range("I126").Value = 7150119
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 21450357
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -21447887
$ws.Range("N126").Value = -28940

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H60").Value = 5143.75
$ws.Range("I60").Value = 256.25
$ws.Range("J60").Value = 7587.5
$ws.Range("K60").Value = 768.75
$ws.Range("L60").Value = 22762.5
$ws.Range("M60").Value = -517.75
$ws.Range("N60").Value = -23264.5
$ws.Range("H107").Value = 173.625
$ws.Range("I107").Value = 173.625
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 520.875
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1399.125
$ws.Range("N107").ClearContents()
$ws.Range("H125").Value = 9133.143
$ws.Range("I125").Value = 7000
$ws.Range("J125").Value = 9239.799999999999
$ws.Range("K125").Value = 21000
$ws.Range("L125").Value = 27719.4
$ws.Range("M125").Value = -16080
$ws.Range("N125").Value = -37559.39999999999
$ws.Range("H131").Value = 577.25
$ws.Range("I131").Value = 324.29413
$ws.Range("J131").Value = 968.1818
$ws.Range("K131").Value = 972.88239
$ws.Range("L131").Value = 2904.5454
$ws.Range("M131").Value = 4067.11761
$ws.Range("N131").Value = -12984.5454
$ws.Range("H140").Value = 1756.6666
$ws.Range("I140").Value = 1330
$ws.Range("J140").Value = 3250
$ws.Range("K140").Value = 3990
$ws.Range("L140").Value = 9750
$ws.Range("M140").Value = 1190
$ws.Range("N140").Value = -20110

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H117").Value = 25787.428
$ws.Range("J117").Value = 25787.428
$ws.Range("L117").Value = 25787.428
$ws.Range("N117").Value = -32671.428

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 1640.5758
$ws.Range("I7").Value = 1288.9333
$ws.Range("J7").Value = 1933.6111
$ws.Range("K7").Value = 1288.9333
$ws.Range("L7").Value = 1933.6111
$ws.Range("M7").Value = -1176.9333
$ws.Range("N7").Value = -2157.6111
$ws.Range("H40").Value = 1981.25
$ws.Range("I40").Value = 1640
$ws.Range("J40").Value = 3005
$ws.Range("K40").Value = 1640
$ws.Range("L40").Value = 3005
$ws.Range("M40").Value = -1504
$ws.Range("N40").Value = -3277
$ws.Range("H46").Value = 3299.75
$ws.Range("I46").Value = 5499.5
$ws.Range("J46").Value = 1100
$ws.Range("K46").Value = 5499.5
$ws.Range("L46").Value = 1100
$ws.Range("M46").Value = -5311.5
$ws.Range("N46").Value = -1476
$ws.Range("H61").Value = 2519.2104
$ws.Range("I61").Value = 2260
$ws.Range("J61").Value = 3901.6667
$ws.Range("K61").Value = 2260
$ws.Range("L61").Value = 3901.6667
$ws.Range("M61").Value = -2058
$ws.Range("N61").Value = -4305.6667
$ws.Range("H113").Value = 2519.2104
$ws.Range("I113").Value = 2260
$ws.Range("J113").Value = 3901.6667
$ws.Range("K113").Value = 2260
$ws.Range("L113").Value = 3901.6667
$ws.Range("M113").Value = -90
$ws.Range("N113").Value = -8241.6667
$ws.Range("H126").Value = 1640.5758
$ws.Range("I126").Value = 1288.9333
$ws.Range("J126").Value = 1933.6111
$ws.Range("K126").Value = 3866.7999
$ws.Range("L126").Value = 5800.8333
$ws.Range("M126").Value = -1396.7999
$ws.Range("N126").Value = -10740.8333
$ws.Range("H132").Value = 18528988
$ws.Range("I132").Value = 38479210
$ws.Range("J132").Value = 3782
$ws.Range("K132").Value = 115437630
$ws.Range("L132").Value = 11346
$ws.Range("M132").Value = -115435100
$ws.Range("N132").Value = -16406

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H107").Value = 1310.0454
$ws.Range("I107").Value = 988.625
$ws.Range("J107").Value = 2167.1667
$ws.Range("K107").Value = 2965.875
$ws.Range("L107").Value = 6501.500100000001
$ws.Range("M107").Value = -1045.875
$ws.Range("N107").Value = -10341.5001
$ws.Range("H113").Value = 432.8
$ws.Range("J113").Value = 683.8333
$ws.Range("L113").Value = 2051.4999
$ws.Range("N113").Value = -6391.4999
$ws.Range("H122").Value = 1639.2727
$ws.Range("I122").Value = 1559.2222
$ws.Range("J122").Value = 1999.5
$ws.Range("K122").Value = 4677.6666
$ws.Range("L122").Value = 5998.5
$ws.Range("M122").Value = -2227.6666
$ws.Range("N122").Value = -10898.5
$ws.Range("H126").Value = 4877.8887
$ws.Range("I126").Value = 5337
$ws.Range("J126").Value = 1205
$ws.Range("K126").Value = 16011
$ws.Range("L126").Value = 3615
$ws.Range("M126").Value = -13541
$ws.Range("N126").Value = -8555
$ws.Range("H132").Value = 1372.5103
$ws.Range("I132").Value = 755.7778
$ws.Range("J132").Value = 3080.3845
$ws.Range("K132").Value = 2267.3334
$ws.Range("L132").Value = 9241.1535
$ws.Range("M132").Value = 262.6666
$ws.Range("N132").Value = -14301.1535
